$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.489800810813904
$ws.Range("B1").Value = 3.749367952346802
$ws.Range("C1").Value = 1.783222913742065
$ws.Range("D1").Value = 1.191227674484253
$ws.Range("E1").Value = 0.7522534132003784
